$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E (which hold text-formatted numbers/percentages as
# inline strings) to Text format first so Excel doesn't auto-convert
# numeric-looking strings (e.g. '0.9998', '244.65') into real numbers.
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '29.575.67'
$ws.Range('E2').Value = '  +2.27%  '
$ws.Range('D3').Value = '1.858.22'
$ws.Range('E3').Value = '  +1.43%  '
$ws.Range('D4').Value = '0.9998'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '244.65'
$ws.Range('E5').Value = '  +1.42%  '
$ws.Range('D6').Value = '0.6938'
$ws.Range('E6').Value = '  +1.01%  '
$ws.Range('D7').Value = '1.001'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '0.07690'
$ws.Range('E8').Value = '  +0.42%  '
$ws.Range('D9').Value = '0.3057'
$ws.Range('E9').Value = '  +0.23%  '
$ws.Range('D10').Value = '23.67'
$ws.Range('E10').Value = '  +0.45%  '
$ws.Range('D11').Value = '0.07762'
$ws.Range('E11').Value = '  -0.57%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.885.94'
$ws.Range('E12').Value = '  +2.76%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '5.151'
$ws.Range('E13').Value = '  +1.40%  '
$ws.Range('D14').Value = '91.58'
$ws.Range('E14').Value = '  +1.26%  '
$ws.Range('D15').Value = '0.6909'
$ws.Range('E15').Value = '  +1.95%  '
$ws.Range('D16').Value = '6.564'
$ws.Range('E16').Value = '  +1.59%  '
$ws.Range('D17').Value = '29.595.25'
$ws.Range('E17').Value = '  +2.30%  '
$ws.Range('D18').Value = '0.000008286'
$ws.Range('E18').Value = '  -0.07%  '
$ws.Range('D19').Value = '2.110.65'
$ws.Range('E19').Value = '  +1.58%  '
$ws.Range('D20').Value = '239.92'
$ws.Range('E20').Value = '  -1.18%  '
$ws.Range('D21').Value = '12.76'
$ws.Range('E21').Value = '  +0.65%  '
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('D23').Value = '7.585'
$ws.Range('E23').Value = '  +2.03%  '
$ws.Range('D24').Value = '1.001'
$ws.Range('E24').Value = '  +0.14%  '
$ws.Range('D25').Value = '0.1495'
$ws.Range('E25').Value = '  +1.71%  '
$ws.Range('D26').Value = '8.912'
$ws.Range('E26').Value = '  +1.40%  '
$ws.Range('D27').Value = '159.82'
$ws.Range('E27').Value = '  -1.03%  '
$ws.Range('D28').Value = '18.26'
$ws.Range('E28').Value = '  +0.38%  '
$ws.Range('E29').Value = '  +0.17%  '
$ws.Range('D30').Value = '4.248'
$ws.Range('E30').Value = '  +0.78%  '
$ws.Range('D31').Value = '4.180'
$ws.Range('E31').Value = '  +1.52%  '
$ws.Range('D32').Value = '1.202'
$ws.Range('E32').Value = '  +1.17%  '
$ws.Range('D33').Value = '0.05094'
$ws.Range('E33').Value = '  -0.53%  '
$ws.Range('D34').Value = '0.7714'
$ws.Range('E34').Value = '  +2.48%  '
$ws.Range('D35').Value = '1.890'
$ws.Range('E35').Value = '  +3.20%  '
$ws.Range('D36').Value = '1.150'
$ws.Range('E36').Value = '  +0.44%  '
$ws.Range('D37').Value = '2.687'
$ws.Range('E37').Value = '  +0.47%  '
$ws.Range('D38').Value = '1.329.20'
$ws.Range('E38').Value = '  +7.75%  '
$ws.Range('D39').Value = '0.01870'
$ws.Range('E39').Value = '  +1.21%  '
$ws.Range('E40').Value = '  +1.11%  '
$ws.Range('D41').Value = '0.9636'
$ws.Range('E41').Value = '  +4.70%  '
$ws.Range('D42').Value = '106.40'
$ws.Range('E42').Value = '  -1.86%  '
$ws.Range('D43').Value = '5.785'
$ws.Range('E43').Value = '  +2.49%  '
$ws.Range('D44').Value = '1.002'
$ws.Range('E44').Value = '  +0.22%  '
$ws.Range('E45').Value = '  +3.86%  '
$ws.Range('D46').Value = '9.769'
$ws.Range('E46').Value = '  +2.76%  '
$ws.Range('D47').Value = '1.999.05'
$ws.Range('E47').Value = '  +1.08%  '
$ws.Range('D48').Value = '0.5218'
$ws.Range('E48').Value = '  +0.91%  '
$ws.Range('D49').Value = '1.771'
$ws.Range('E49').Value = '  +2.03%  '
$ws.Range('D50').Value = '63.53'
$ws.Range('E50').Value = '  -1.09%  '
$ws.Range('D51').Value = '6.948'
$ws.Range('E51').Value = '  +0.73%  '

# Restore the default 'Normal' style so no stray number-format style
# index is left attached to these cells (matches original formatting).
$ws.Range('D2:E51').Style = 'Normal'
